# Update timestamp column (O) from "2022-12-25 12:56:16" to "2022-12-25 20:51:25"
# for all data rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldTimestamp = "2022-12-25 12:56:16"
$newTimestamp = "2022-12-25 20:51:25"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)  # column O = 15
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
